$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D:D").Insert()
Write-Host "done"
